# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / handoff / handback datetime
# stamps for the e1f29d49 (and related) rows across the Overview, zh-cn and
# de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
# Row 3 corresponds to e1f29d49-765f-4abe-8a5d-7a268dab63cd.md
$wsOverview.Range("G3").Value = "2016-09-06 07:49:30"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Row 2 corresponds to a7560759-331f-4243-83c1-b7e49ded50f9.md
$wsZhCn.Range("K2").Value = "2016-09-06 07:50:04"
# Row 3 corresponds to e1f29d49-765f-4abe-8a5d-7a268dab63cd.md
$wsZhCn.Range("H3").Value = "2016-09-06 07:49:19"
$wsZhCn.Range("K3").Value = "2016-09-06 07:50:04"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
# Row 2 corresponds to a7560759-331f-4243-83c1-b7e49ded50f9.md
$wsDeDe.Range("K2").Value = "2016-09-06 07:50:33"
# Row 3 corresponds to e1f29d49-765f-4abe-8a5d-7a268dab63cd.md
$wsDeDe.Range("H3").Value = "2016-09-06 07:49:30"
$wsDeDe.Range("K3").Value = "2016-09-06 07:50:33"
